# ---------------------------------------------------------------------------
# Applies the "Sync attendance_reports ... 2026-01-12 11:15:25" edit to the
# session-analysis workbook:
#   1. Swap "<email>, System" -> "System, <email>" in every "Recorded By"
#      (column G) cell that has that exact text.
#   2. Update the summary statistics block (L6, L8, L9, L10).
#   3. Update the per-group statistics rows for B1-10/11/12 seminar rows
#      (O/Q/R/S on rows 16-18 and 24-26).
#   4. Six sessions (B1-7..B1-12, session 12) moved from "Pending" to
#      "Recorded": fill color, Recorded-By, Students and Status all change.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap "email, System" -> "System, email" in column G -----------------
$swapped = 0
for ($r = 1; $r -le 260; $r++) {
    $cell = $ws.Range("G" + $r)
    if ($cell.Text -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
        $swapped = $swapped + 1
    }
}
Write-Host "Swapped Recorded-By order on" $swapped "rows"

# --- 2. Summary statistics block --------------------------------------------
$ws.Range("L6").Value = 138
$ws.Range("L8").Value = 48
$ws.Range("L9").Value = "'53.5%"
$ws.Range("L10").Value = "'79.3%"

# --- 3. Per-group statistics for B1-10, B1-11, B1-12 (rows 16-18) ----------
# and B1-7, B1-8, B1-9 (rows 24-26)
$groupRows = @(16, 17, 18, 24, 25, 26)
$sValues = @{
    16 = "'73.7%"
    17 = "'63.9%"
    18 = "'82.1%"
    24 = "'71.3%"
    25 = "'75.9%"
    26 = "'74.7%"
}
foreach ($r in $groupRows) {
    $ws.Range("O" + $r).Value = 12
    $ws.Range("Q" + $r).Value = 3
    $ws.Range("R" + $r).Value = "'57.1%"
    $ws.Range("S" + $r).Value = $sValues[$r]
}

# --- 4. Pending -> Recorded sessions ----------------------------------------
$recordedRows = @{
    33  = "27/31"
    54  = "16/18"
    75  = "16/21"
    196 = "19/27"
    217 = "24/29"
    238 = "25/29"
}

foreach ($r in $recordedRows.Keys) {
    # Copy the green "Recorded" formatting (A2:I2) onto A{r}:I{r}, replacing
    # the yellow "Pending" formatting, without touching the existing values.
    $ws.Range("A2:I2").Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Range("G" + $r).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $r).Value = $recordedRows[$r]
    $ws.Range("I" + $r).Value = "Recorded"
}

$excel.CutCopyMode = 0

Write-Host "Done applying edits"
